# Junction_Flooding_126 edit:
#  - round row 5 data (B5:AH5) to 2 decimal places ("custom accuracy")
#  - delete row 6 (trim to "1000 data points" -> fewer sample rows here)
#  - tighten several column widths (consequence of the narrower numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in row 5 (columns B..AH) to 2 decimal places.
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $current = [double]$cell.Value2
    $cell.Value = [Math]::Round($current, 2)
}

# Delete row 6 entirely (shifts any rows below it up, none here).
$ws.Rows.Item(6).Delete()

# Adjust column widths to match the new (narrower) content widths.
# NOTE: Excel quantizes ColumnWidth (characters) to the nearest pixel using
# the workbook's Maximum Digit Width before storing the OOXML <col width>,
# so asking for exactly "7" characters round-trips to ~7.83. Nudging the
# requested width slightly past the integer boundary (x.1) lands on the
# clean integer pixel-width bucket that serializes as a whole number.
$narrowCols = @(2,3,7,9,10,11,12,13,15,16,17,22,23,24,26,27,28,29,30,33,34)
foreach ($c in $narrowCols) {
    $ws.Columns.Item($c).ColumnWidth = 6.1
}
$ws.Columns.Item(20).ColumnWidth = 7.1
